$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the two "isOwner" validation strings so they reflect the owner's own
# view of their newly-created IPA containers (was incorrectly false).
$ws.Range("J3").Value = "status=200||type=ipa_ss||name=First IPA Container by Project Neon1||desc=First SSE - IPA Container created by postman||userid=(SYS_USER2)||ispublic=false||isOwner=true"
$ws.Range("J4").Value = "status=200||type=ipa_ss||name=Second IPA Container by Project Neon1||desc=Second SSE - IPA Container created by postman||userid=(SYS_USER2)||ispublic=false||isOwner=true"

# Remove the obsolete OPQA-4635 "last viewed time" test row (was row 19);
# everything below shifts up by one.
$ws.Rows.Item(19).Delete()

# Keep the sheet's saved selection in sync with the now-shorter data range.
$ws.Range("L2:L25").Select() | Out-Null
